# Resolutions de plusieurs problemes sur la metaheuristique
#
# 1. Rename the "Notre solution" header (column G) to "Notre Heuristique"
# 2. Add a new column I "Notre MetaHeuristique" with its own result values
# 3. Update the existing column G values (heuristic results changed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Step 1: rename existing header text ---
$ws.Range("G2").Value = "Notre Heuristique"

# --- Step 2: new header for column I, matching the bold header style used
#     by the rest of row 2 (copy format from G2, which already carries it) ---
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I2").Value = "Notre MetaHeuristique"

# --- Step 3: updated "Notre Heuristique" (column G) results ---
$gValues = @{
    3  = 444
    4  = 786
    5  = 732
    6  = 300
    7  = 756
    8  = 522
    9  = 396
    10 = 528
    11 = 348
    12 = 726
    13 = 300
    14 = 840
    15 = 160
    16 = 241
    17 = 229
    18 = 378
    19 = 309
    20 = 414
    21 = 404
    22 = 398
    23 = 345
    24 = 206
    25 = 642
    26 = 585
    27 = 624
    28 = 543
    29 = 631
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# --- Step 4: new "Notre MetaHeuristique" (column I) results ---
$iValues = @{
    3  = 336
    4  = 324
    5  = 288
    6  = 222
    7  = 390
    8  = 438
    9  = 330
    10 = 396
    11 = 354
    12 = 540
    13 = 294
    14 = 750
    15 = 128
    16 = 142
    17 = 181
    18 = 343
    19 = 253
    20 = 305
    21 = 303
    22 = 351
    23 = 233
    24 = 229
    25 = 486
    26 = 425
    27 = 553
    28 = 443
    29 = 572
}

foreach ($row in $iValues.Keys) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
}

# --- Step 5: adjust the view so the last edited cell is selected ---
$ws.Range("E33").Select() | Out-Null
